# Update reviews for parisk
# - Blank out the "polite_expressions" value (C13) that used to read "nan".
# - Append a new review row (row 14) with a fresh annotation entry.

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# C13 becomes an empty text value (was "nan").
# A leading apostrophe forces Excel to store it as text rather than clearing
# the cell outright, matching the other blank "polite_expressions" cells
# elsewhere in the column (e.g. C2, C3, ...).
$ws.Range("C13").Value = "'"

# New row 14: another parisk review entry.
$ws.Range("A14").Value = "parisk"
# politeness_score is stored as text "3" (not a number) for this row.
$ws.Range("B14").Value = "'3"
$ws.Range("C14").Value = "nan"
$ws.Range("D14").Value = "DIS"
$ws.Range("E14").Value = "WRI"
$ws.Range("F14").Value = "3a6bf25f-9f71-48b7-a40b-7e968e5f9337"
$ws.Range("G14").Value = "ry-TW-WAb_annotated.xlsx"
$ws.Range("H14").Value = "I suggest to change it to e.g. 'from the true to the approximate posterior' to avoid confusion."
$ws.Range("I14").Value = "Correct"
